$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.462.07'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +2.76%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.605.43'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.65'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.521'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +6.70%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '26.82'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +6.39%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +2.71%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0601'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +2.59%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0910'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.88%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.834.39'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +2.53%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.614.13'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.94%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '29.449.75'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.534'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +3.43%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.70'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.82%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '63.44'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +3.59%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '240.71'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +5.19%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +3.58%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +1.66%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.09'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '154.44'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.86%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +4.84%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.25'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +3.14%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +2.47%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.22'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.62%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.412.91'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +1.38%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.97%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.83'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +5.10%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +2.64%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +3.77%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.99'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +1.91%  '
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +5.74%  '
$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '53.20'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +22.97%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.796'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +3.32%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '65.89'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +2.91%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.27'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.743.87'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +2.70%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.858'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '86.65'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.82%  '
